# "Generate Report for Handoff"
#
# A new handoff report run updated the "ht" (handoff-type) priority flag
# and the "Latest Handoff/HO Xliff Generate" timestamps for the rows that
# previously reported an error (the ones whose priority did not match the
# handoff type). These values are mirrored across the Overview sheet and
# the per-locale (zh-cn / de-de) sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Rows (matching the source files that hit the handback-priority mismatch
# error) whose Priority/Handoff-date values were refreshed by the new run.
$rows = @(8, 10, 11, 12, 13, 14)

$zhcnTimestamp = "2016-08-18 18:22:19"
$dedeTimestamp = "2016-08-18 18:22:25"

foreach ($r in $rows) {
    # Priority column (E) now reports "ht" instead of being blank.
    $zhcn.Range("E$r").Value = "ht"
    $dede.Range("E$r").Value = "ht"

    # Latest Handoff Datetime (H) on each locale sheet.
    $zhcn.Range("H$r").Value = $zhcnTimestamp
    $dede.Range("H$r").Value = $dedeTimestamp

    # Latest HO Xliff Generate Date (G) on the Overview sheet mirrors the
    # de-de locale's handoff timestamp for these rows.
    $overview.Range("G$r").Value = $dedeTimestamp
}
